# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot run).
# Values that are valid Excel numbers are prefixed with a leading apostrophe
# so they stay text cells (matching the source data's inline-string typing)
# instead of being auto-converted to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.919.18"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.171.83"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'572.31"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'166.19"
$ws.Range("E6").Value = "  -4.06%  "
$ws.Range("D7").Value = "'0.594"
$ws.Range("E7").Value = "  -5.57%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "3.716.69"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "64.779.46"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'25.52"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "3.164.77"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "'414.04"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'12.69"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'7.15"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'68.48"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "'0.201"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  -5.32%  "
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'1.83"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").Value = "'21.36"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'4.94"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'6.33"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "'1.13"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").Value = "'154.99"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").Value = "2.715.61"
$ws.Range("E36").Value = "  -3.61%  "
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("E38").Value = "  -6.32%  "
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "'5.57"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'291.83"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("D45").Value = "'21.35"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'0.0986"
$ws.Range("E48").Value = "  -9.97%  "
$ws.Range("D49").Value = "'10.45"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'0.899"
$ws.Range("E51").Value = "  -3.87%  "
